$d = $word.ActiveDocument

# 1. Rename the "Approach" heading to "Methods".
$d.Content.Find.Execute("Approach", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Methods", 2) | Out-Null

# 2. Word tracks the location of the last text edit with a hidden
#    "_GoBack" bookmark. Move it from its old spot (end of the
#    "Fig. 1....." paragraph) to right after the text we just typed
#    ("Methods"), collapsed (zero-length), matching what real Word does
#    when you type over a selection.
#
#    A truly collapsed range built exactly at "end-of-paragraph-text"
#    cannot be passed straight into Bookmarks.Add here, so instead we
#    insert a temporary placeholder character, wrap the bookmark around
#    it, then delete the placeholder -- the bookmark collapses in place,
#    exactly where we need it.
$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.TrimEnd([char]13) -eq "Methods") {
        $headingPara = $p
        break
    }
}

$target = $headingPara.Range.Duplicate
$target.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$target.Collapse(0) | Out-Null      # 0 = wdCollapseEnd -> right after "Methods"
$target.InsertAfter("x") | Out-Null
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
$goBack = $d.Bookmarks("_GoBack")
$goBack.Range.Text = ""

# 3. The footer's cached PAGE field result is re-rendered to "2" (the
#    doc now spills onto a second page once "Approach" becomes "Methods"
#    plus the bookmark housekeeping above). Update the cached field text
#    directly -- it lives inside a page-number content control (SDT).
foreach ($sec in $d.Sections) {
    $footer = $sec.Footers(1)
    if ($footer.Exists) {
        $fieldResult = $footer.Range.Paragraphs(1).Range.Words(1)
        if ($fieldResult.Text.Trim() -eq "1") {
            $fieldResult.Text = "2"
        }
    }
}
